$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the
# 27c80a73-16fb-4437-a628-5ab6f9ace938.md row (row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-17 15:52:53"

# zh-cn sheet: update Correspond Handoff Datetime (H2) and
# Correspond Handback DateTime (K2) for the 27c80a73 row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-17 15:52:31"
$wsZhCn.Range("K2").Value = "2016-10-17 15:53:36"

# de-de sheet: update Correspond Handoff Datetime (H2) and
# Correspond Handback DateTime (K2) for the 27c80a73 row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-17 15:52:53"
$wsDeDe.Range("K2").Value = "2016-10-17 15:54:14"
